$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("展览").Range("F8").Value = 2717
$wb.Worksheets.Item("展览").Range("F11").Value = 2404
$wb.Worksheets.Item("展览").Range("F18").Value = 6815
$wb.Worksheets.Item("展览").Range("F20").Value = 89
$wb.Worksheets.Item("展览").Range("F24").Value = 7692
$wb.Worksheets.Item("展览").Range("F38").Value = 2575
$wb.Worksheets.Item("展览").Range("F42").Value = 1148
$wb.Worksheets.Item("展览").Range("F45").Value = 3608
$wb.Worksheets.Item("展览").Range("F47").Value = 1157
$wb.Worksheets.Item("展览").Range("F48").Value = 96
$wb.Worksheets.Item("演出").Range("F5").Value = 234
$wb.Worksheets.Item("全部类型").Range("F7").Value = 2717
$wb.Worksheets.Item("全部类型").Range("F9").Value = 234
$wb.Worksheets.Item("全部类型").Range("F11").Value = 2404
$wb.Worksheets.Item("全部类型").Range("F20").Value = 6815
$wb.Worksheets.Item("全部类型").Range("F22").Value = 89
$wb.Worksheets.Item("全部类型").Range("F25").Value = 7692
$wb.Worksheets.Item("全部类型").Range("F38").Value = 2575
$wb.Worksheets.Item("全部类型").Range("F41").Value = 1148
$wb.Worksheets.Item("全部类型").Range("F45").Value = 3608
$wb.Worksheets.Item("全部类型").Range("F48").Value = 1157
$wb.Worksheets.Item("全部类型").Range("F49").Value = 96
